$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "Resolving-Mac" sending-cluster rows (previously rows 14-17);
# remaining rows shift up to close the gap.
$ws.Range("A14:T17").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# New TPM-derived values for the remaining rows (now rows 2-13)
$arr = New-Object 'object[,]' 12,20
# row 2: Sending=ECs, Target=ECs
$arr[0,0] = 'ECs'
$arr[0,1] = 'Lgi2'
$arr[0,2] = 'Adam23'
$arr[0,3] = 'ECs'
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 0.126156
$arr[0,7] = 0.378468
$arr[0,8] = 0.01239214403266014
$arr[0,9] = 0.01239214403266014
$arr[0,10] = 1
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.2189473333333334
$arr[0,13] = 0.656842
$arr[0,14] = 0.009402596261870986
$arr[0,15] = 0.009402596261870984
$arr[0,16] = 0.02762151978400001
$arr[0,17] = 0.248593678056
$arr[0,18] = 0.0001165183271580571
$arr[0,19] = 0.000116518327158057
# row 3: Sending=ECs, Target=FAPs
$arr[1,0] = 'ECs'
$arr[1,1] = 'Lgi2'
$arr[1,2] = 'Adam23'
$arr[1,3] = 'FAPs'
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 0.126156
$arr[1,7] = 0.378468
$arr[1,8] = 0.01239214403266014
$arr[1,9] = 0.01239214403266014
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 20.07911633333334
$arr[1,13] = 60.237349
$arr[1,14] = 0.8622887582286424
$arr[1,15] = 0.8622887582286423
$arr[1,16] = 2.533101000148001
$arr[1,17] = 22.797909001332
$arr[1,18] = 0.01068560648971299
$arr[1,19] = 0.01068560648971299
# row 4: Sending=ECs, Target=MuSCs
$arr[2,0] = 'ECs'
$arr[2,1] = 'Lgi2'
$arr[2,2] = 'Adam23'
$arr[2,3] = 'MuSCs'
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.126156
$arr[2,7] = 0.378468
$arr[2,8] = 0.01239214403266014
$arr[2,9] = 0.01239214403266014
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 2.823530666666667
$arr[2,13] = 8.470592
$arr[2,14] = 0.1212552739852724
$arr[2,15] = 0.1212552739852723
$arr[2,16] = 0.3562053347840001
$arr[2,17] = 3.205848013056
$arr[2,18] = 0.001502612819945163
$arr[2,19] = 0.001502612819945163
# row 5: Sending=ECs, Target=Resolving-Mac
$arr[3,0] = 'ECs'
$arr[3,1] = 'Lgi2'
$arr[3,2] = 'Adam23'
$arr[3,3] = 'Resolving-Mac'
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 0.126156
$arr[3,7] = 0.378468
$arr[3,8] = 0.01239214403266014
$arr[3,9] = 0.01239214403266014
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 0.1642436666666667
$arr[3,13] = 0.492731
$arr[3,14] = 0.007053371524214274
$arr[3,15] = 0.007053371524214274
$arr[3,16] = 0.020720324012
$arr[3,17] = 0.186482916108
$arr[3,18] = 0.00008740639584392686
$arr[3,19] = 0.00008740639584392685
# row 6: Sending=FAPs, Target=ECs
$arr[4,0] = 'FAPs'
$arr[4,1] = 'Lgi2'
$arr[4,2] = 'Adam23'
$arr[4,3] = 'ECs'
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 9.690151
$arr[4,7] = 29.070453
$arr[4,8] = 0.9518512547181717
$arr[4,9] = 0.9518512547181717
$arr[4,10] = 1
$arr[4,11] = 0.3333333333333333
$arr[4,12] = 0.2189473333333334
$arr[4,13] = 0.656842
$arr[4,14] = 0.009402596261870986
$arr[4,15] = 0.009402596261870984
$arr[4,16] = 2.121632721047334
$arr[4,17] = 19.094694489426
$arr[4,18] = 0.008949873049470288
$arr[4,19] = 0.008949873049470286
# row 7: Sending=FAPs, Target=FAPs
$arr[5,0] = 'FAPs'
$arr[5,1] = 'Lgi2'
$arr[5,2] = 'Adam23'
$arr[5,3] = 'FAPs'
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 9.690151
$arr[5,7] = 29.070453
$arr[5,8] = 0.9518512547181717
$arr[5,9] = 0.9518512547181717
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 20.07911633333334
$arr[5,13] = 60.237349
$arr[5,14] = 0.8622887582286424
$arr[5,15] = 0.8622887582286423
$arr[5,16] = 194.5696692165664
$arr[5,17] = 1751.127022949097
$arr[5,18] = 0.8207706364493074
$arr[5,19] = 0.8207706364493073
# row 8: Sending=FAPs, Target=MuSCs
$arr[6,0] = 'FAPs'
$arr[6,1] = 'Lgi2'
$arr[6,2] = 'Adam23'
$arr[6,3] = 'MuSCs'
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 9.690151
$arr[6,7] = 29.070453
$arr[6,8] = 0.9518512547181717
$arr[6,9] = 0.9518512547181717
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 2.823530666666667
$arr[6,13] = 8.470592
$arr[6,14] = 0.1212552739852724
$arr[6,15] = 0.1212552739852723
$arr[6,16] = 27.36043851313067
$arr[6,17] = 246.243946618176
$arr[6,18] = 0.1154169846840772
$arr[6,19] = 0.1154169846840772
# row 9: Sending=FAPs, Target=Resolving-Mac
$arr[7,0] = 'FAPs'
$arr[7,1] = 'Lgi2'
$arr[7,2] = 'Adam23'
$arr[7,3] = 'Resolving-Mac'
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 9.690151
$arr[7,7] = 29.070453
$arr[7,8] = 0.9518512547181717
$arr[7,9] = 0.9518512547181717
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 0.1642436666666667
$arr[7,13] = 0.492731
$arr[7,14] = 0.007053371524214274
$arr[7,15] = 0.007053371524214274
$arr[7,16] = 1.591545930793667
$arr[7,17] = 14.323913377143
$arr[7,18] = 0.00671376053531678
$arr[7,19] = 0.00671376053531678
# row 10: Sending=MuSCs, Target=ECs
$arr[8,0] = 'MuSCs'
$arr[8,1] = 'Lgi2'
$arr[8,2] = 'Adam23'
$arr[8,3] = 'ECs'
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 0.3640136666666667
$arr[8,7] = 1.092041
$arr[8,8] = 0.03575660124916825
$arr[8,9] = 0.03575660124916825
$arr[8,10] = 1
$arr[8,11] = 0.3333333333333333
$arr[8,12] = 0.2189473333333334
$arr[8,13] = 0.656842
$arr[8,14] = 0.009402596261870986
$arr[8,15] = 0.009402596261870984
$arr[8,16] = 0.07969982161355557
$arr[8,17] = 0.7172983945220001
$arr[8,18] = 0.0003362048852426408
$arr[8,19] = 0.0003362048852426407
# row 11: Sending=MuSCs, Target=FAPs
$arr[9,0] = 'MuSCs'
$arr[9,1] = 'Lgi2'
$arr[9,2] = 'Adam23'
$arr[9,3] = 'FAPs'
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 0.3640136666666667
$arr[9,7] = 1.092041
$arr[9,8] = 0.03575660124916825
$arr[9,9] = 0.03575660124916825
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 20.07911633333334
$arr[9,13] = 60.237349
$arr[9,14] = 0.8622887582286424
$arr[9,15] = 0.8622887582286423
$arr[9,16] = 7.309072759923223
$arr[9,17] = 65.781654839309
$arr[9,18] = 0.03083251528962201
$arr[9,19] = 0.03083251528962201
# row 12: Sending=MuSCs, Target=MuSCs
$arr[10,0] = 'MuSCs'
$arr[10,1] = 'Lgi2'
$arr[10,2] = 'Adam23'
$arr[10,3] = 'MuSCs'
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 0.3640136666666667
$arr[10,7] = 1.092041
$arr[10,8] = 0.03575660124916825
$arr[10,9] = 0.03575660124916825
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 2.823530666666667
$arr[10,13] = 8.470592
$arr[10,14] = 0.1212552739852724
$arr[10,15] = 0.1212552739852723
$arr[10,16] = 1.027803750919111
$arr[10,17] = 9.250233758272
$arr[10,18] = 0.004335676481250028
$arr[10,19] = 0.004335676481250027
# row 13: Sending=MuSCs, Target=Resolving-Mac
$arr[11,0] = 'MuSCs'
$arr[11,1] = 'Lgi2'
$arr[11,2] = 'Adam23'
$arr[11,3] = 'Resolving-Mac'
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 0.3640136666666667
$arr[11,7] = 1.092041
$arr[11,8] = 0.03575660124916825
$arr[11,9] = 0.03575660124916825
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 0.1642436666666667
$arr[11,13] = 0.492731
$arr[11,14] = 0.007053371524214274
$arr[11,15] = 0.007053371524214274
$arr[11,16] = 0.05978693933011112
$arr[11,17] = 0.5380824539710001
$arr[11,18] = 0.0002522045930535679
$arr[11,19] = 0.0002522045930535679

$ws.Range("A2:T13").Value = $arr
